$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename existing sheets and add the two new ones in the right order
# ---------------------------------------------------------------------
$wsWydatki = $wb.Worksheets.Item(1)
$wsWydatki.Name = "Wydatki"

$wsManex = $wb.Worksheets.Item(2)

$wsTauron = $wb.Worksheets.Item(3)
$wsTauron.Name = "Rachunki Tauron"

$wsZgk = $wb.Worksheets.Add($null, $wsTauron)
$wsZgk.Name = "Rachunki ZGK"

$wsKaty = $wb.Worksheets.Add($null, $wsZgk)
$wsKaty.Name = "Podatek Kąty"

# ---------------------------------------------------------------------
# 2. "Rachunki Tauron" sheet (sheet3) - template table for Tauron bills
#    Cell writes are ordered to match the original shared-string order.
# ---------------------------------------------------------------------
$wsTauron.Range("B2").Value = "Data płatności"
$wsTauron.Range("C2").Value = "Nr faktury"
$wsTauron.Range("D2").Value = "Kwota"
$wsTauron.Range("E2").Value = "Zapłacono"

$wsTauron.Range("B3").Value = 41442
$wsTauron.Range("B3").NumberFormat = "mm-dd-yy"
$wsTauron.Range("C3").Value = "559021187/3/s"
$wsTauron.Range("D3").Value = 58.23

$wsTauron.Range("B4").Value = 41470
$wsTauron.Range("B3").Copy()
$wsTauron.Range("B4").PasteSpecial(-4122)
$wsTauron.Range("C4").Value = "559021187/4/s"
$wsTauron.Range("D4").Value = 58.23

$wsTauron.Range("F2").Value = "Kwota przelewu"

$loTauron = $wsTauron.ListObjects.Add(1, $wsTauron.Range("B2:F14"), $null, 1)
$loTauron.Name = "Tabela2"
$loTauron.TableStyle = "TableStyleMedium7"

$wsTauron.Columns.Item(2).ColumnWidth = 16.75
$wsTauron.Columns.Item(3).ColumnWidth = 19.75
$wsTauron.Columns.Item(4).ColumnWidth = 14.5
$wsTauron.Columns.Item(5).ColumnWidth = 17.75
$wsTauron.Columns.Item(6).ColumnWidth = 16.75

$wsTauron.Range("C19").Select()

# ---------------------------------------------------------------------
# 3. "Rachunki ZGK" sheet (sheet4) - template table for ZGK bills
# ---------------------------------------------------------------------
$wsZgk.Range("B2").Value = "Data płatności"
$wsZgk.Range("C2").Value = "Kwota"
$wsZgk.Range("D2").Value = "Opis"
$wsZgk.Range("E2").Value = "Nr faktury"
$wsZgk.Range("F2").Value = "Zapłacono"
$wsZgk.Range("G2").Value = "Kwota przelewu"

$loZgk = $wsZgk.ListObjects.Add(1, $wsZgk.Range("B2:G12"), $null, 1)
$loZgk.Name = "Tabela3"
$loZgk.TableStyle = "TableStyleMedium5"

$wsZgk.Columns.Item(2).ColumnWidth = 14.75
$wsZgk.Columns.Item(3).ColumnWidth = 14.5
$wsZgk.Columns.Item(4).ColumnWidth = 16.75
$wsZgk.Columns.Item(5).ColumnWidth = 14.25
$wsZgk.Columns.Item(6).ColumnWidth = 14.75
$wsZgk.Columns.Item(7).ColumnWidth = 16.75

$wsZgk.Rows.Item(2).Select()

# ---------------------------------------------------------------------
# 4. "Podatek Kąty" sheet (sheet5) - Kąty tax payment schedule
# ---------------------------------------------------------------------
$wsKaty.Range("B2").Value = "Termin zapłaty"
$wsKaty.Range("C2").Value = "Kwota"
$wsKaty.Range("D2").Value = "Zapłacono"

$wsKaty.Range("B3").Value = 41348
$wsTauron.Range("B3").Copy()
$wsKaty.Range("B3").PasteSpecial(-4122)
$wsKaty.Range("C3").Value = 11

$wsKaty.Range("B4").Value = 41409
$wsKaty.Range("B4").PasteSpecial(-4122)
$wsKaty.Range("C4").Value = 9

$wsKaty.Range("B5").Value = 41532
$wsKaty.Range("B5").PasteSpecial(-4122)
$wsKaty.Range("C5").Value = 9

$wsKaty.Range("B6").Value = 41593
$wsKaty.Range("B6").PasteSpecial(-4122)
$wsKaty.Range("C6").Value = 9

$loKaty = $wsKaty.ListObjects.Add(1, $wsKaty.Range("B2:D9"), $null, 1)
$loKaty.Name = "Tabela46"
$loKaty.TableStyle = "TableStyleMedium3"

$wsKaty.Columns.Item(2).ColumnWidth = 30.75
$wsKaty.Columns.Item(3).ColumnWidth = 16
$wsKaty.Columns.Item(4).ColumnWidth = 36.25

$wsKaty.Range("B15").Select()

# ---------------------------------------------------------------------
# 5. "Faktury Manex" sheet loses the tabSelected flag (handled above as a
#    side effect of moving the active sheet/tab to "Podatek Kąty").
# ---------------------------------------------------------------------
